# Update the bulk-upload template header row so that the "roleIds" column
# becomes "roleKeys" to support the new dynamic member role set up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: A1=email (unchanged), B1 was "roleIds" -> now "roleKeys",
# C1 stays "remark".
$ws.Range("B1").Value = "roleKeys"
$ws.Range("C1").Value = "remark"

# Move the active selection to C6 (matches the saved selection in the file).
$ws.Range("C6").Select()
